$d = $word.ActiveDocument

# Locate the paragraph that still ends with "LOB1018: Física I (Requisito fraco)"
# and delete everything from the start of the next paragraph (the blank line
# right after it) through the end of the "© 2020 ..." paragraph, i.e. the
# blank paragraph, the "Ver no Jupiter..." paragraph and the "© 2020..."
# paragraph.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -match "LOB1018: F.sica I \(Requisito fraco\)") {
        $startPara = $i + 1
        $found = $true
        break
    }
}

if ($found) {
    $endPara = $startPara + 2
    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
